# Auto-generated script applying scheduled-runner market data refresh
# Updates currentAveragePrice* / Leve Price* / Leve Profit* columns (H-N)
# for the rows whose market snapshot changed, across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 340
$ws.Range("I33").Value = 387.5
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 387.5
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = -158.5
$ws.Range("N33").Value = -608

$ws.Range("H40").Value = 3948
$ws.Range("I40").Value = 6299.3335
$ws.Range("K40").Value = 6299.3335
$ws.Range("M40").Value = -6124.3335

$ws.Range("H49").Value = 15000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H51").Value = 8596
$ws.Range("I51").Value = 8604.571
$ws.Range("J51").Value = 8566
$ws.Range("K51").Value = 8604.571
$ws.Range("L51").Value = 8566
$ws.Range("M51").Value = -8120.571
$ws.Range("N51").Value = -9534

$ws.Range("H88").Value = 11112361
$ws.Range("I88").Value = 50000400
$ws.Range("J88").Value = 1493.1428
$ws.Range("K88").Value = 50000400
$ws.Range("L88").Value = 1493.1428
$ws.Range("M88").Value = -49999994
$ws.Range("N88").Value = -2305.1428

$ws.Range("H91").Value = 11112361
$ws.Range("I91").Value = 50000400
$ws.Range("J91").Value = 1493.1428
$ws.Range("K91").Value = 50000400
$ws.Range("L91").Value = 1493.1428
$ws.Range("M91").Value = -49998996
$ws.Range("N91").Value = -4301.1428

$ws.Range("H96").Value = 757.4
$ws.Range("I96").Value = 695.6667
$ws.Range("J96").Value = 850
$ws.Range("K96").Value = 2087.0001
$ws.Range("L96").Value = 2550
$ws.Range("M96").Value = -714.0001000000002
$ws.Range("N96").Value = -5296

$ws.Range("H111").Value = 1692.7188
$ws.Range("I111").Value = 473.6
$ws.Range("J111").Value = 1918.4814
$ws.Range("K111").Value = 1420.8
$ws.Range("L111").Value = 5755.4442
$ws.Range("M111").Value = 1646.2
$ws.Range("N111").Value = -11889.4442

$ws.Range("H113").Value = 3365.6365
$ws.Range("I113").Value = 3336.6667
$ws.Range("J113").Value = 3376.5
$ws.Range("K113").Value = 3336.6667
$ws.Range("L113").Value = 3376.5
$ws.Range("M113").Value = -82.66670000000022
$ws.Range("N113").Value = -9884.5

$ws.Range("H116").Value = 8010.2915
$ws.Range("I116").Value = 2349.6
$ws.Range("J116").Value = 9499.947
$ws.Range("K116").Value = 2349.6
$ws.Range("L116").Value = 9499.947
$ws.Range("M116").Value = 1092.4
$ws.Range("N116").Value = -16383.947

$ws.Range("H121").Value = 2112.7144
$ws.Range("J121").Value = 2112.7144
$ws.Range("L121").Value = 6338.1432
$ws.Range("N121").Value = -9832.143199999999

$ws.Range("H137").Value = 10021321
$ws.Range("I137").Value = 16667200
$ws.Range("K137").Value = 50001600
$ws.Range("M137").Value = -49999050

$ws.Range("H141").Value = 6394.3125
$ws.Range("I141").Value = 3022.1428
$ws.Range("K141").Value = 9066.428400000001
$ws.Range("M141").Value = -3886.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 970
$ws.Range("I2").Value = 926.1818
$ws.Range("K2").Value = 926.1818
$ws.Range("M2").Value = -813.1818

$ws.Range("H39").Value = 30000
$ws.Range("I39").Value = 30000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 30000
$ws.Range("L39").ClearContents()
$ws.Range("M39").Value = -29480
$ws.Range("N39").Value = 0

$ws.Range("H61").Value = 4088.6667
$ws.Range("I61").Value = 2709.9167
$ws.Range("K61").Value = 2709.9167
$ws.Range("M61").Value = -2497.9167

$ws.Range("H74").Value = 280696.34
$ws.Range("I74").Value = 619068.25
$ws.Range("K74").Value = 619068.25
$ws.Range("M74").Value = -618194.25

$ws.Range("H77").Value = 280696.34
$ws.Range("I77").Value = 619068.25
$ws.Range("K77").Value = 3095341.25
$ws.Range("M77").Value = -3090973.25

$ws.Range("H102").Value = 5412.231
$ws.Range("I102").Value = 5666
$ws.Range("K102").Value = 5666
$ws.Range("M102").Value = -4044

$ws.Range("H116").Value = 970
$ws.Range("I116").Value = 926.1818
$ws.Range("K116").Value = 926.1818
$ws.Range("M116").Value = 1367.8182

$ws.Range("H125").Value = 126332.664
$ws.Range("J125").Value = 126332.664
$ws.Range("L125").Value = 126332.664
$ws.Range("N125").Value = -136172.664

$ws.Range("H132").Value = 2573.5833
$ws.Range("I132").Value = 1792.6111
$ws.Range("K132").Value = 5377.8333
$ws.Range("M132").Value = -2847.8333

$ws.Range("H136").Value = 4088.6667
$ws.Range("I136").Value = 2709.9167
$ws.Range("K136").Value = 8129.750100000001
$ws.Range("M136").Value = -5579.750100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 970
$ws.Range("I3").Value = 926.1818
$ws.Range("K3").Value = 926.1818
$ws.Range("M3").Value = -812.1818

$ws.Range("H107").Value = 1561.8235
$ws.Range("I107").Value = 1428.3572
$ws.Range("J107").Value = 2184.6667
$ws.Range("K107").Value = 1428.3572
$ws.Range("L107").Value = 2184.6667
$ws.Range("M107").Value = 491.6428000000001
$ws.Range("N107").Value = -6024.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4460.3335
$ws.Range("I31").Value = 3343.8
$ws.Range("J31").Value = 5856
$ws.Range("K31").Value = 3343.8
$ws.Range("L31").Value = 5856
$ws.Range("M31").Value = -3048.8
$ws.Range("N31").Value = -6446

$ws.Range("H34").Value = 4460.3335
$ws.Range("I34").Value = 3343.8
$ws.Range("J34").Value = 5856
$ws.Range("K34").Value = 3343.8
$ws.Range("L34").Value = 5856
$ws.Range("M34").Value = -3141.8
$ws.Range("N34").Value = -6260

$ws.Range("H105").Value = 2173
$ws.Range("I105").Value = 1562.5
$ws.Range("J105").Value = 2347.4285
$ws.Range("K105").Value = 1562.5
$ws.Range("L105").Value = 2347.4285
$ws.Range("M105").Value = 184.5
$ws.Range("N105").Value = -5841.4285

$ws.Range("H107").Value = 866.86664
$ws.Range("I107").Value = 539.4
$ws.Range("J107").Value = 1521.8
$ws.Range("K107").Value = 539.4
$ws.Range("L107").Value = 1521.8
$ws.Range("M107").Value = 1380.6
$ws.Range("N107").Value = -5361.8

$ws.Range("H122").Value = 4452.3076
$ws.Range("I122").Value = 3278.1428
$ws.Range("J122").Value = 5822.1665
$ws.Range("K122").Value = 9834.428400000001
$ws.Range("L122").Value = 17466.4995
$ws.Range("M122").Value = -7384.428400000001
$ws.Range("N122").Value = -22366.4995

$ws.Range("H132").Value = 3569.0715
$ws.Range("I132").Value = 3360.4285
$ws.Range("J132").Value = 3777.7144
$ws.Range("K132").Value = 10081.2855
$ws.Range("L132").Value = 11333.1432
$ws.Range("M132").Value = -7551.2855
$ws.Range("N132").Value = -16393.1432

$ws.Range("H134").Value = 1938.46
$ws.Range("I134").Value = 1857.9773
$ws.Range("J134").Value = 2528.6667
$ws.Range("K134").Value = 5573.9319
$ws.Range("L134").Value = 7586.000100000001
$ws.Range("M134").Value = -3038.9319
$ws.Range("N134").Value = -12656.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws.Range("H112").Value = 7341.6665
$ws.Range("J112").Value = 7499
$ws.Range("L112").Value = 22497
$ws.Range("N112").Value = -24713

$ws.Range("H131").Value = 23629.908
$ws.Range("J131").Value = 2066.3333
$ws.Range("L131").Value = 6198.999899999999
$ws.Range("N131").Value = -16278.9999

$ws.Range("H134").Value = 2512.1538
$ws.Range("I134").Value = 1471.5
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 4414.5
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = 655.5
$ws.Range("N134").Value = -55140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3548.6
$ws.Range("I113").Value = 3200
$ws.Range("K113").Value = 3200
$ws.Range("M113").Value = -1030

$ws.Range("H122").Value = 9048.857
$ws.Range("I122").Value = 7081.3335
$ws.Range("J122").Value = 10524.5
$ws.Range("K122").Value = 21244.0005
$ws.Range("L122").Value = 31573.5
$ws.Range("M122").Value = -18794.0005
$ws.Range("N122").Value = -36473.5

$ws.Range("H132").Value = 3433.3333
$ws.Range("I132").Value = 1650
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 4950
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -2420
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2763.1667
$ws.Range("I7").Value = 2763.1667
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2763.1667
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2651.1667

$ws.Range("H40").Value = 25604.738
$ws.Range("I40").Value = 31146.756
$ws.Range("K40").Value = 31146.756
$ws.Range("M40").Value = -31010.756

$ws.Range("H46").Value = 1830.1
$ws.Range("I46").Value = 1001
$ws.Range("J46").Value = 1922.2222
$ws.Range("K46").Value = 1001
$ws.Range("L46").Value = 1922.2222
$ws.Range("M46").Value = -813
$ws.Range("N46").Value = -2298.2222

$ws.Range("H126").Value = 2763.1667
$ws.Range("I126").Value = 2763.1667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8289.500100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -5819.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4468.2354
$ws.Range("I81").Value = 1610
$ws.Range("K81").Value = 3220
$ws.Range("M81").Value = -2159

$ws.Range("H84").Value = 4468.2354
$ws.Range("I84").Value = 1610
$ws.Range("K84").Value = 16100
$ws.Range("M84").Value = -10796

$ws.Range("H126").Value = 1021.25
$ws.Range("I126").Value = 1021.25
$ws.Range("K126").Value = 3063.75
$ws.Range("M126").Value = -593.75

$ws.Range("H136").Value = 76928430
$ws.Range("I136").Value = 100001360
$ws.Range("K136").Value = 300004080
$ws.Range("M136").Value = -300001530
